# Regenerate s_val data to filter save games.
# Updates the numeric stat columns (B:E) and the derived sum column (G)
# for rows 2-6 on the active worksheet. Column F (Win) is left unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.6545652718822623, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 3.536033448013082)
    3 = @(1.445647641019636, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 3.755628166162433)
    4 = @(0.6545652718822623, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 3.536033448013082)
    5 = @(0.6545652718822623, 1.626987699542094, 18.71679738969934, 13.86384647080068, 34.86219683192438)
    6 = @(1.445647641019636, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 4.327115817150455)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
    $ws.Range("E$row").Value = $vals[3]
    $ws.Range("G$row").Value = $vals[4]
}
